$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix the tiny floating point re-computation on row 4's timestamp ---
$ws.Cells.Item(4, 1).Value2 = 45729.48582644676

# --- Row 5: new history entry (single "calcium" deficiency) ---
$ws.Cells.Item(5, 1).Value2 = 45729.49783483797
$ws.Cells.Item(5, 1).NumberFormat = $ws.Cells.Item(4, 1).NumberFormat
$ws.Cells.Item(5, 2).Value = "Aarti"
$ws.Cells.Item(5, 3).Value = 25
$ws.Cells.Item(5, 4).Value = "Male"
$ws.Cells.Item(5, 5).Value = 50
$ws.Cells.Item(5, 6).Value = 1.5
$ws.Cells.Item(5, 7).Value = 22.22
$ws.Cells.Item(5, 8).Value = "Normal weight - Maintain a balanced diet and exercise."
$ws.Cells.Item(5, 9).Value = "Veg"
$ws.Cells.Item(5, 10).Value = "calcium"

$k5 = @"

Recommendations for calcium Deficiency:
Tofu, raw, firm, prepared with calcium sulfate
Cheese, Mexican, blend, reduced fat
Cheese, cheddar, nonfat or fat free
Cheese, swiss, low fat
Cheese, swiss, low sodium
Cheese, mozzarella, part skim milk
Cheese, gruyere
Cheese, monterey
Cheese, port de salut
Cheese, swiss
Cheese, swiss
Cheese, provolone, sliced
Cheese, provolone, reduced fat
Cheese, monterey jack, solid
Cheese, low-sodium, cheddar or colby
Cheese, muenster
Cheese, mozzarella, low sodium
Cheese, provolone
Cheese, monterey, low fat
Cheese, brick
Cheese, mexican, queso asadero
Cheese, colby
Cheese, Mexican blend
Cheese, Swiss, nonfat or fat free
Cheese, queso fresco, solid
Cheese, cheddar
Cheese, mexican, queso chihuahua
Cheese, cheddar, sharp, sliced
Cheese, cheddar
Cheese, white, queso blanco
Cheese, mozzarella, nonfat
Cheese, cheddar, reduced fat
Cheese, tilsit
Cheese, parmesan, grated, refrigerated
Cheese, cheshire
Cheese, parmesan, hard
Cheese, caraway
Imitation cheese, american or cheddar, low cholesterol
Cheese, fontina
Cheese, mexican, queso anejo
"@
$ws.Cells.Item(5, 11).Value = $k5

# --- Row 6: new history entry (two deficiencies: calcium, vitamin_E) ---
$ws.Cells.Item(6, 1).Value2 = 45729.51109761735
$ws.Cells.Item(6, 1).NumberFormat = $ws.Cells.Item(4, 1).NumberFormat
$ws.Cells.Item(6, 2).Value = "Aarti"
$ws.Cells.Item(6, 3).Value = 25
$ws.Cells.Item(6, 4).Value = "Male"
$ws.Cells.Item(6, 5).Value = 50
$ws.Cells.Item(6, 6).Value = 1.5
$ws.Cells.Item(6, 7).Value = 22.22
$ws.Cells.Item(6, 8).Value = "Normal weight - Maintain a balanced diet and exercise."
$ws.Cells.Item(6, 9).Value = "Veg"
$ws.Cells.Item(6, 10).Value = "calcium, vitamin_E"

$k6 = @"

Recommendations for calcium and vitamin_E Deficiency:
Cheese, swiss
Cheese, swiss
Cheese, cheddar, sharp, sliced
Cheese, Mexican, blend, reduced fat
Grape leaves, raw
Cheese, cheddar
Cheese, swiss, low sodium
Cheese, cheddar
Tofu, raw, firm, prepared with calcium sulfate
Cheese, gruyere
Cheese, monterey
Turnip greens, raw
Cheese, port de salut
Cheese, mozzarella, part skim milk
Cheese, cheddar, nonfat or fat free
Cheese, white, queso blanco
Cheese, swiss, low fat
Cheese, cheddar, reduced fat
Cheese, low-sodium, cheddar or colby
Cheese, muenster
Peppers, jalapeno, raw
Cheese, provolone
Cheese, colby
Cheese, brick
Cheese, provolone, reduced fat
Cheese, mexican, queso asadero
Cheese, Mexican blend
Cheese, monterey, low fat
Cheese, mexican, queso chihuahua
Cheese, mozzarella, low sodium
Cheese, mozzarella, low moisture, part-skim, shredded
Cheese, mozzarella, low moisture, part-skim
Cheese food, pasteurized process, American, vitamin D fortified
Cheese, pasteurized process, American, without added vitamin D
Cheese, provolone, sliced
Cheese, mozzarella, low moisture, part-skim
Cheese product, pasteurized process, American, vitamin D fortified
Cheese, monterey jack, solid
"@
$ws.Cells.Item(6, 11).Value = $k6
